$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New cell: H6 ---
$ws1.Range("H6").Value = "Ater offset"

# --- Row 7: update the "Test -20cm" note + measured triple ---
$ws1.Range("F7").Value = "Test -20cm z, Left Back"
$ws1.Range("G7").Value = "339.5, 234.5, 70.9 "

# --- Row 22 (Calibrating L Back Hip) : replace the long note in B22 ---
$ws1.Range("B22").Value = "90 deg on servo angle is the axis 0 mark.0 deg is upright. 120 is down.  So 339 input should yield a 21 deg angle below the horizontal. So servo should receive angle of 111. Offset: (360 - theta1) + 90 deg"

# --- Row 23 (Calibrating L back knee) : replace the long note in B23 ---
$ws1.Range("B23").Value = "actual knee only really has a bit over 90 deg of range of motion. et servo 120 point to be when knee is all the way at the bottom, perpendicular to ground. 234 deg should yield 36 deg above this point, Offset: 120 -  (270 - theta 2)  "

# --- Row 24 (Calibrating L Back Ankle) : replace B24 and append new commentary columns C24:F24 ---
$ws1.Range("B24").Value = "Cant do -0.22 elevationg cause of ROM concerne with ankle so im going to try using 20cm which states an angle of 67 degrees. If I set the bottoming point to be 120 deg on servo (point where ankle is max range of motion when knee is at 90deg)  then for this leg, 160 servo angle is where ankle and knee are at 90 deg, which correponds to the models theta3 = 90. So if i want a theta 3 of 67 deg. if theta 3 is less than 90, i do 160 - (90 - theta 3). if theta 3 is greater than 90 i do: 160 - theta 3"
$ws1.Range("C24").Value = "Goddamn, if max range on code is 180, but max on servo is 120 I need there to be a ratio. Y = (x/120)* 180"
$ws1.Range("D24").Value = "ankle set to 120 deg for bottomed out position at max range. Min angle is around 70deg. Models 90 deg point is 100on the servo. So if I want a theta of 71 deg.  I could do 71 + 10"
$ws1.Range("E24").Value = "the thing I need to consider is that if the knee angle changes, the ankle still maintains its angle relative to the floot, and also the range of motion changes. if I change the theta 2 angle by 70, then I need to adjust the end effector theta 3 angle by 70 just to maintain the position. so i coul dsay when the model tells theta 2 to change by an angle, i add that offset to theta 3. theta 2- theta2 '  offet to 3? might be a more elgean tway to implemetn this tho"
$ws1.Range("F24").Value = "I want the middle of the mechanical range and the middle of the servo range to line up. "

# New cells C24:F24 need the same wrapped style as the rest of column B in that block
$ws1.Range("C24:F24").WrapText = $true

# --- Column widths for the two newly-used columns E & F ---
$ws1.Columns.Item(5).ColumnWidth = 38.5
$ws1.Columns.Item(6).ColumnWidth = 23.5

# --- Row heights stay the same as authored (re-assert defensively) ---
$ws1.Rows.Item(22).RowHeight = 137.4
$ws1.Rows.Item(23).RowHeight = 216.6
$ws1.Rows.Item(24).RowHeight = 224.4

# --- Make Sheet1 the active / selected sheet (was "Ankle Calibration") ---
$ws1.Activate()
$ws1.Range("F24").Select()
